$d = $word.ActiveDocument

function Replace-Exact($findText, $replaceText) {
    $range = $d.Content
    $range.Find.ClearFormatting()
    $range.Find.Execute($findText, $true, $true, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null
}

# 1. Merge "SE CONSIDERAN TODAS LAS BAJAS COMO LOGICAS =>" + " TODAS LAS TABLAS..." into one run
Replace-Exact "SE CONSIDERAN TODAS LAS BAJAS COMO LOGICAS => TODAS LAS TABLAS POSEEN UNA COLUMNA ACTIVO QUE SE PONE EN FALSE CUANDO SE LO ELIMINA." "SE CONSIDERAN TODAS LAS BAJAS COMO LOGICAS => TODAS LAS TABLAS POSEEN UNA COLUMNA ACTIVO QUE SE PONE EN FALSE CUANDO SE LO ELIMINA."

# 2. Merge "EN LA CONSULTA DE UPDATE" + " PUBLICACION, NO SE LE PASARON..." into one run
Replace-Exact "EN LA CONSULTA DE UPDATE PUBLICACION, NO SE LE PASARON LOS PARAMETROS QUE SABEMOS POR EL NEGOCIO QUE NO PUEDEN VARIAR (ID_USUARIO, FECHA_INICIO,  FECHA_VENCIMIENTO)" "EN LA CONSULTA DE UPDATE PUBLICACION, NO SE LE PASARON LOS PARAMETROS QUE SABEMOS POR EL NEGOCIO QUE NO PUEDEN VARIAR (ID_USUARIO, FECHA_INICIO,  FECHA_VENCIMIENTO)"

# 3. Merge "EN LOS FILTROS DE BUSQUEDA...CASO CONTRARIO NO," + " NO PERMITIMOS HACER CONSULTAS LIKE..." into one run
Replace-Exact "EN LOS FILTROS DE BUSQUEDA, CUANDO ESTA ULTIMA ES EXACTA, PERMITE FILTRAR POR CAMPOS NUMERICOS, (CASO CONTRARIO NO, NO PERMITIMOS HACER CONSULTAS LIKE CON CAMPOS NUMERICOS)" "EN LOS FILTROS DE BUSQUEDA, CUANDO ESTA ULTIMA ES EXACTA, PERMITE FILTRAR POR CAMPOS NUMERICOS, (CASO CONTRARIO NO, NO PERMITIMOS HACER CONSULTAS LIKE CON CAMPOS NUMERICOS)"

# 4. Merge "NO HACEMOS NINGUNA ELIMINACION FISICA DE " + "LOS REGISTROS (DELETE QUERIES)..." into one run
Replace-Exact "NO HACEMOS NINGUNA ELIMINACION FISICA DE LOS REGISTROS (DELETE QUERIES) SINO QUE MARCAMOS COMO ACTIVO FALSE EL REGISTRO EN LA TABLA CORRESPONDIENTE" "NO HACEMOS NINGUNA ELIMINACION FISICA DE LOS REGISTROS (DELETE QUERIES) SINO QUE MARCAMOS COMO ACTIVO FALSE EL REGISTRO EN LA TABLA CORRESPONDIENTE"

# 5. Add new paragraph at the end about migration of qualifications
$last = $d.Paragraphs.Last
$last.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$newPara.Range.Text = "EN LA MIGRACIÓN TOMAMOS A LAS CALIFICACIONES Y LAS DIVIDIMOS POR DOS PARA ADAPTARLAS A LA NUEVA VERSIÓN DE LAS CALIFICACIONES (QUE VAN DE 1 A 5)"
